# daily auto push: 2026-01-09 18:47 UTC
# Insert two new rows before row 601 so the following rows (previously
# 601-642) shift down to 603-644, then populate the four "new" rows
# (601-604) with their correct data. Everything from row 605 onward keeps
# its original values automatically (just shifted by the insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 601:642 down by two rows, creating two blank rows at 601:602.
$ws.Rows("601:602").Insert()

# Helper-free direct writes for the four rows that now hold new data.
# Column A holds date strings that look like dates ("2026/01/09" etc.);
# force them to stay as literal text (matching the rest of the sheet)
# by writing under a Text number format, then clearing the format so no
# style index is left behind on the cell.

$ws.Range("A601").NumberFormat = "@"
$ws.Range("A601").Value = "2026/01/09"
$ws.Range("A601").ClearFormats()
$ws.Range("B601").Value = "金"
$ws.Range("C601").Value = 23
$ws.Range("D601").Value = 201

$ws.Range("A602").NumberFormat = "@"
$ws.Range("A602").Value = "2026/01/10"
$ws.Range("A602").ClearFormats()
$ws.Range("B602").Value = "土"
$ws.Range("C602").Value = 2
$ws.Range("D602").Value = 201

$ws.Range("A603").NumberFormat = "@"
$ws.Range("A603").Value = "2026/01/10"
$ws.Range("A603").ClearFormats()
$ws.Range("B603").Value = "土"
$ws.Range("C603").Value = 13
$ws.Range("D603").Value = 201

$ws.Range("A604").NumberFormat = "@"
$ws.Range("A604").Value = "2026/01/10"
$ws.Range("A604").ClearFormats()
$ws.Range("B604").Value = "土"
$ws.Range("C604").Value = 16
$ws.Range("D604").Value = 201
